# Updates the cryptos worksheet with refreshed Price and Volume(1h) figures
# (GitHub Actions symbol-list refresh, 2023-01-18 13:30 UTC).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume columns store plain text in this sheet (e.g. "301.60", "-0.08%").
# Mark each target cell as Text first so Excel keeps the new values - including
# trailing zeros like "301.20" or percent signs like "-0.23%" - as literal text
# instead of silently re-parsing them into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"

# Apply the refreshed values.
$ws.Range("D2").Value = "301.20"
$ws.Range("E2").Value = "-0.23%"
$ws.Range("D3").Value = "32.45"
$ws.Range("E3").Value = "1.98%"
$ws.Range("D4").Value = "5.069"
$ws.Range("E4").Value = "-1.36%"
$ws.Range("D5").Value = "0.07673"
$ws.Range("E5").Value = "-2.17%"
$ws.Range("D6").Value = "2.113"
$ws.Range("E6").Value = "-6.34%"
$ws.Range("D7").Value = "7.848"
$ws.Range("E7").Value = "0.41%"
$ws.Range("D8").Value = "3.779"
$ws.Range("E8").Value = "-0.82%"
$ws.Range("D9").Value = "0.9207"
$ws.Range("E9").Value = "-0.77%"
$ws.Range("E10").Value = "-0.46%"
$ws.Range("D11").Value = "0.07944"
$ws.Range("E11").Value = "4.02%"
$ws.Range("D12").Value = "0.08402"
$ws.Range("E12").Value = "-5.19%"
$ws.Range("D13").Value = "0.03059"
$ws.Range("E13").Value = "-1.42%"
$ws.Range("D14").Value = "0.09986"
$ws.Range("E14").Value = "-0.30%"
$ws.Range("D15").Value = "0.001516"
$ws.Range("E15").Value = "0.08%"
$ws.Range("D16").Value = "0.005866"
$ws.Range("E16").Value = "1.29%"
$ws.Range("D18").Value = "3.462"
$ws.Range("E18").Value = "-0.59%"
$ws.Range("E19").Value = "-4.47%"
$ws.Range("D20").Value = "0.3341"
$ws.Range("E20").Value = "1.49%"
$ws.Range("D21").Value = "0.1326"
$ws.Range("E21").Value = "-0.14%"
$ws.Range("D22").Value = "4.290"
$ws.Range("E22").Value = "-0.63%"
$ws.Range("D23").Value = "0.1976"
$ws.Range("E23").Value = "10.34%"
$ws.Range("E24").Value = "-1.40%"
$ws.Range("D25").Value = "0.001238"
$ws.Range("E25").Value = "-1.02%"
$ws.Range("D26").Value = "0.004125"
$ws.Range("E26").Value = "-7.86%"
$ws.Range("D27").Value = "0.0001251"
$ws.Range("E27").Value = "0.12%"
$ws.Range("D39").Value = "0.01710"
$ws.Range("E39").Value = "-3.98%"
$ws.Range("D40").Value = "0.04672"
$ws.Range("E40").Value = "-2.62%"
$ws.Range("D41").Value = "0.007449"
$ws.Range("E41").Value = "0.87%"
$ws.Range("D42").Value = "0.1353"
$ws.Range("E42").Value = "-0.70%"
$ws.Range("D43").Value = "0.002331"
$ws.Range("E43").Value = "6.52%"
$ws.Range("D44").Value = "0.01053"
$ws.Range("E44").Value = "7.22%"
$ws.Range("D45").Value = "0.00006208"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").Value = "0.12%"
$ws.Range("D47").Value = "0.003001"
$ws.Range("E47").Value = "-62.42%"
$ws.Range("D48").Value = "0.8233"
$ws.Range("E48").Value = "17.10%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.12%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.12%"
